$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Data.NetCall" -> "Data.NetCalls" (row 10, Containing Type column for SendServerInfo)
$ws.Cells.Item(10, 3).Value2 = "Data.NetCalls"

# Insert a new blank row above the old "MulticastKitUpdated" row (worksheet row 103),
# pushing it and everything below it down by one row.
$ws.Rows.Item(103).Insert()

# The sheet's Excel table ("Table2") needs to grow to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E110"))

# Row 103: first new ListSync.NetCalls entry (replaces the old KitSync.NetCalls
# "MulticastKitUpdated" call with the start of the new list-sync implementation).
$ws.Cells.Item(103, 1).Value2 = "MulticastListItemUpdated"
$ws.Cells.Item(103, 2).Value2 = 3000
$ws.Cells.Item(103, 3).Value2 = "ListSync.NetCalls"
$ws.Cells.Item(103, 4).Value2 = "FROM_EITHER"
$ws.Cells.Item(103, 5).Value2 = "ushort syncId, int pk"

# Row 104: second new ListSync.NetCalls entry.
$ws.Cells.Item(104, 1).Value2 = "MulticastListItemsUpdated"
$ws.Cells.Item(104, 2).Value2 = 3001
$ws.Cells.Item(104, 3).Value2 = "ListSync.NetCalls"
$ws.Cells.Item(104, 4).Value2 = "FROM_EITHER"
$ws.Cells.Item(104, 5).Value2 = "ushort syncId, int[] pks"

# Reflect the final scroll position / selection from the author's editing session.
$excel.Goto($ws.Range("A82"), $true)
$ws.Range("E93").Select()
